$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue {
    param($range, [string]$text)
    # Force the value to be stored as text (not auto-coerced to a number),
    # while leaving the cell's style untouched.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

# Row 2: NLP engineer job (id stays "3", score 100 -> 97, reason text updated)
Set-TextValue $ws.Range("A2") "3"
$ws.Range("B2").Value = 97
$ws.Range("C2").Value = "The job as an NLP engineer involves analyzing and preprocessing large-scale text data, and conducting PoCs in areas such as LLM, NLP, DL, ML, and object detection/classification. Your experience with projects like LLMGuard (bias detection), and Multi Model Data Analysis (data analysis) makes you well-suited for this role. Your skills in NLP, Pytorch, Computer Vision, and Python align closely with the required skills for this job, hence the high score of 97."

# Row 3: Frontend Engineer Intern (id "1" -> "2", score 75 -> 71, reason text updated)
Set-TextValue $ws.Range("A3") "2"
$ws.Range("B3").Value = 71
$ws.Range("C3").Value = "The job of Frontend Engineer Intern is suitable for you because your experience in ReactJS, JavaScript, and CSS aligns well with the required skills for the job. Additionally, your projects involving ReactJS development make you a suitable candidate for ensuring efficient and visually appealing web design and user experience. The moderate score suggests that while you meet the basic requirements, there may be other factors influencing the hiring decision."

# Row 4: SDE Intern (id "2" -> "1", score 62 -> 71, reason text updated)
Set-TextValue $ws.Range("A4") "1"
$ws.Range("B4").Value = 71
$ws.Range("C4").Value = "The job of SDE Intern has a score of 71. This indicates that it is moderately suitable for you. The job requires skills such as MongoDB, ReactJS, JavaScript, Web Development, and NodeJS, which align with your experience in projects like LLMGuard and Literature Society IITJ Website. However, the job may not be a perfect fit as it may require additional skills or experience in certain areas."
